# Excel COM-interop edit script
# Commit message: "cleaning up commits... changing guest pricing"
#
# The Variables sheet holds a GUEST_MARKUP constant (cell A5) that is used
# to calculate each product's guest sales price (ffcsaGuestSalesPrice,
# column M on the Pricelist sheet) from its purchase price. This commit
# bumps GUEST_MARKUP from 0.9 to 0.92 and refreshes the guest sales price
# for every product whose price actually changes as a result (rows whose
# guest price was already 0 are left untouched, matching the source diff).

$wb = $excel.ActiveWorkbook

# --- 1. Update the GUEST_MARKUP constant on the "Variables" sheet ---
$wsVars = $wb.Worksheets.Item("Variables")
$wsVars.Range("A5").Value = 0.92

# --- 2. Recalculated ffcsaGuestSalesPrice values (column M) on "Pricelist" ---
# Map of row number -> new guest sales price, taken from the recalculation
# that results from the GUEST_MARKUP change above.
$guestPriceUpdates = @{
    2 = 15.07
    3 = 15.07
    4 = 14.03
    5 = 15.59
    6 = 17.15
    7 = 17.15
    8 = 20.78
    9 = 22.44
    10 = 17.66
    11 = 12.47
    12 = 10.39
    13 = 23.85
    14 = 27.35
    15 = 48.58
    16 = 63.83
    17 = 48.53
    19 = 37.67
    20 = 15.82
    21 = 27.02
    22 = 49.36
    23 = 17.66
    24 = 25.98
    25 = 14.55
    26 = 33.77
    27 = 56.53
    28 = 11.43
    29 = 24.32
    30 = 29.61
    31 = 24.32
    32 = 29.72
    33 = 15.59
    34 = 22.29
    35 = 45.59
    36 = 34.45
    38 = 23.07
    39 = 24.57
    40 = 27.15
    41 = 44.73
    42 = 74.82
    43 = 61.74
    44 = 49.47
    45 = 43.23
    46 = 33.77
    47 = 16.11
    48 = 14.86
    49 = 26.34
    50 = 12.47
    51 = 26.91
    52 = 12.47
    53 = 39.98
    54 = 27.02
    55 = 16.07
    56 = 14.86
    57 = 37.67
    58 = 14.81
    59 = 9.87
    60 = 9.87
    61 = 12.47
    62 = 28.06
    63 = 20.13
    64 = 21.82
    65 = 51.44
    66 = 45.41
    67 = 38.45
    68 = 31.04
    69 = 66.77
    70 = 56.17
    71 = 43.64
    72 = 19.22
    73 = 12.47
    74 = 20.78
    75 = 10.39
    76 = 22.44
    77 = 43.64
    78 = 9.38
    79 = 27.37
    80 = 20.94
    81 = 31.72
    82 = 24.68
    83 = 44.42
    84 = 51.96
    86 = 22.34
    87 = 62.55
    88 = 15.07
    89 = 18.18
    90 = 20.26
    91 = 16.11
    92 = 20.26
    93 = 20.26
    94 = 22.34
    95 = 22.34
    96 = 17.15
    97 = 20.78
    98 = 22.34
    99 = 18.7
    100 = 18.7
    101 = 20.78
    102 = 18.7
    103 = 20.78
    104 = 9.35
    105 = 20.26
    106 = 20.78
    107 = 17.15
    108 = 11.43
    109 = 17.66
    110 = 20.26
    111 = 18.7
    112 = 12.47
    113 = 7.27
    114 = 10.39
    115 = 20.78
    116 = 17.15
    117 = 17.15
    118 = 17.15
    119 = 17.66
    120 = 17.15
    121 = 20.78
    122 = 20.78
    123 = 20.78
    124 = 15.59
    125 = 20.78
    126 = 20.78
    127 = 24.94
    128 = 20.78
    129 = 10.39
    130 = 15.59
    131 = 10.39
    132 = 16.63
    133 = 10.39
    134 = 22.86
    135 = 20.78
    136 = 17.25
    137 = 21.58
    138 = 22.74
    139 = 19.98
    140 = 22.74
    141 = 17.25
    142 = 21.58
    143 = 22.75
    144 = 17.25
    145 = 19.98
    146 = 22.74
    147 = 15.98
    148 = 9.59
    149 = 14.39
    150 = 23.55
    151 = 19.98
    152 = 21.58
    153 = 30.78
    154 = 33.97
    155 = 23.98
    156 = 22.38
    157 = 17.58
    158 = 15.98
    159 = 15.98
    160 = 15.98
    161 = 9.27
    162 = 19.18
    163 = 11.19
    164 = 19.18
}

$wsPrice = $wb.Worksheets.Item("Pricelist")

foreach ($row in $guestPriceUpdates.Keys) {
    $wsPrice.Cells.Item($row, 13).Value = $guestPriceUpdates[$row]
}

Write-Host "Updated GUEST_MARKUP to 0.92 and refreshed $($guestPriceUpdates.Count) guest sales price cells."
